$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing block (rows 144-147) down to the new rows (148-151)
$ws.Range("A144:B147").Copy()
$ws.Range("A148:B151").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Set the new values for rows 148-151
$ws.Range("A148").Value = "(prefixo)_DGE_SQE_B3_P9_17_insere_numero_outros_funcionarios_masc"
$ws.Range("B148").Value = "`${(prefixo)_DGESQE_B3_P9_16_existem_otros_funcionarios}=1"

$ws.Range("A149").Value = "(prefixo)_DGE_SQE_B3_P9_18_insere_numero_outros_funcionarios_fem"
$ws.Range("B149").Value = "`${(prefixo)_DGESQE_B3_P9_16_existem_otros_funcionarios}=1"

$ws.Range("A150").Value = "(prefixo)_DGE_SQE_B3_P9_19_insere_numero_outros_funcionarios_masc_especifique"
$ws.Range("B150").Value = "`${(prefixo)_DGESQE_B3_P9_16_existem_otros_funcionarios}=1 and `${QEPE_(prefixo)_DGESQE_B3_P1_15_tipo_sala_outro_numero}>=1"

$ws.Range("A151").Value = "(prefixo)_DGE_SQE_B3_P9_20_insere_numero_outros_funcionarios_fem_especifique"
$ws.Range("B151").Value = "`${(prefixo)_DGESQE_B3_P9_16_existem_otros_funcionarios}=1 and `${QEPE_(prefixo)_DGESQE_B3_P1_15_tipo_sala_outro_numero}>=1"

# Ensure row height matches the rest of the sheet (15, customHeight)
$ws.Range("A148:B151").RowHeight = 15

# Update selection to match the new active region
$ws.Range("A148:XFD151").Select()
